$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New header cells G1, H1, I1 (bold like the rest of the header row)
$ws.Cells.Item(1, 7).Value = "liberal"
$ws.Cells.Item(1, 8).Value = "kapica"
$ws.Cells.Item(1, 9).Value = "peker"

$ws.Range("G1:I1").Font.Bold = $true

# New data cells on row 6
$ws.Cells.Item(6, 7).Value = 9
$ws.Cells.Item(6, 8).Value = 9
$ws.Cells.Item(6, 9).Value = 8

# Update selection to reflect the new active cell location (I7), mirroring the diff
$ws.Range("I7").Select()
